# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Acelga" (Vega Monumental Concepción)
# right after the existing row 120, pushing the remaining historical rows
# down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 121:122 - everything from the old row 121 onward
# (including the dimension) shifts down by two rows automatically.
$ws.Rows("121:122").Insert()

# New row 121 - Primera quality, "Región Metropolitana"
$ws.Cells.Item(121, 1).Value  = 11
$ws.Cells.Item(121, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(121, 3).Value  = "Bíobío"
$ws.Cells.Item(121, 4).Value  = 44567
$ws.Cells.Item(121, 5).Value  = 8
$ws.Cells.Item(121, 6).Value  = 100112009
$ws.Cells.Item(121, 7).Value  = "Acelga"
$ws.Cells.Item(121, 8).Value  = "Sin especificar"
$ws.Cells.Item(121, 9).Value  = "Primera"
$ws.Cells.Item(121, 10).Value = 200
$ws.Cells.Item(121, 11).Value = 600
$ws.Cells.Item(121, 12).Value = 700
$ws.Cells.Item(121, 13).Value = 650
$ws.Cells.Item(121, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(121, 15).Value = "Región Metropolitana"
$ws.Cells.Item(121, 16).Value = 650
$ws.Cells.Item(121, 17).Value = 1
$ws.Cells.Item(121, 18).Value = "Hortaliza"

# New row 122 - Segunda quality, "Región Metropolitana"
$ws.Cells.Item(122, 1).Value  = 11
$ws.Cells.Item(122, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(122, 3).Value  = "Bíobío"
$ws.Cells.Item(122, 4).Value  = 44567
$ws.Cells.Item(122, 5).Value  = 8
$ws.Cells.Item(122, 6).Value  = 100112009
$ws.Cells.Item(122, 7).Value  = "Acelga"
$ws.Cells.Item(122, 8).Value  = "Sin especificar"
$ws.Cells.Item(122, 9).Value  = "Segunda"
$ws.Cells.Item(122, 10).Value = 100
$ws.Cells.Item(122, 11).Value = 500
$ws.Cells.Item(122, 12).Value = 500
$ws.Cells.Item(122, 13).Value = 500
$ws.Cells.Item(122, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(122, 15).Value = "Región Metropolitana"
$ws.Cells.Item(122, 16).Value = 500
$ws.Cells.Item(122, 17).Value = 1
$ws.Cells.Item(122, 18).Value = "Hortaliza"
